$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2024-02-09 Friday" "2024-02-10 Saturday"

Replace-Text "779×6=" "659×9="
Replace-Text "293×5=" "533×3="
Replace-Text "790×9=" "942×2="
Replace-Text "423×2=" "578×7="
Replace-Text "854×2=" "939×8="
Replace-Text "833×5=" "373×5="
Replace-Text "548×3=" "815×7="
Replace-Text "666×6=" "803×8="
Replace-Text "591×8=" "441×2="
Replace-Text "975×2=" "149×3="
Replace-Text "614×4=" "972×5="
Replace-Text "867×9=" "436×9="
Replace-Text "285×8=" "583×9="
Replace-Text "625×7=" "569×8="
Replace-Text "117×3=" "277×6="
Replace-Text "158×4=" "267×3="
Replace-Text "779×5=" "525×7="
Replace-Text "306×3=" "459×9="
Replace-Text "479×6=" "112×9="
Replace-Text "581×2=" "556×2="
Replace-Text "675×9=" "253×6="
Replace-Text "209×4=" "232×8="
Replace-Text "242×3=" "440×4="
Replace-Text "812×8=" "492×5="
Replace-Text "992×6=" "621×9="
